$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cell values from the latest crypto data refresh.
# Cells that hold numeric-looking text (prices) are forced to Text format
# first so Excel does not silently convert them into real numbers and
# lose formatting such as trailing zeros (e.g. "1.00" vs 1).

$ws.Range("D2").Value = "50.787.54"
$ws.Range("E2").Value = "  -1.01%  "
$ws.Range("D3").Value = "2.924.32"
$ws.Range("E3").Value = "  -1.67%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "375.03"
$ws.Range("E5").Value = "  -1.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.91"
$ws.Range("E6").Value = "  -2.27%  "
$ws.Range("E7").Value = "  -1.02%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.584"
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.92"
$ws.Range("E10").Value = "  -2.55%  "
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0843"
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("D13").Value = "3.383.93"
$ws.Range("E13").Value = "  -1.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.95"
$ws.Range("E14").Value = "  -1.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.56"
$ws.Range("E15").Value = "  +0.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "11.31"
$ws.Range("E16").Value = "  +53.59%  "
$ws.Range("D17").Value = "2.925.03"
$ws.Range("E17").Value = "  -1.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.988"
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("D19").Value = "50.761.36"
$ws.Range("E19").Value = "  -0.93%  "
$ws.Range("E20").Value = "  -6.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.36"
$ws.Range("E21").Value = "  -3.32%  "
$ws.Range("D22").Value = "0.0₃0951"
$ws.Range("E22").Value = "  -0.49%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.79"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.78"
$ws.Range("E24").Value = "  +1.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.16"
$ws.Range("E25").Value = "  +9.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.98"
$ws.Range("E26").Value = "  -1.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.40"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "25.46"
$ws.Range("E29").Value = "  -1.54%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.110"
$ws.Range("E30").Value = "  -6.46%  "
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.162"
$ws.Range("E31").Value = "  -3.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.96"
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.83"
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("E34").Value = "  -0.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "33.05"
$ws.Range("E35").Value = "  -3.60%  "
$ws.Range("E36").Value = "  -3.58%  "
$ws.Range("E38").Value = "  +3.88%  "
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.41"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.81"
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("E42").Value = "  -4.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "119.34"
$ws.Range("E43").Value = "  -2.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.96"
$ws.Range("E44").Value = "  -2.04%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.37"
$ws.Range("E45").Value = "  +3.42%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.04"
$ws.Range("E46").Value = "  -2.00%  "
$ws.Range("E47").Value = "  -1.39%  "
$ws.Range("E48").Value = "  -1.76%  "
$ws.Range("D49").Value = "1.986.01"
$ws.Range("E49").Value = "  -2.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0323"
$ws.Range("E50").Value = "  -2.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.18"
$ws.Range("E51").Value = "  +1.61%  "
